$wb = $excel.ActiveWorkbook

# Sheet: 展览 (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 877
$ws1.Cells.Item(3, 6).Value = 259
$ws1.Cells.Item(4, 6).Value = 20
$ws1.Cells.Item(5, 6).Value = 164
$ws1.Cells.Item(7, 6).Value = 157
$ws1.Cells.Item(8, 6).Value = 4640
$ws1.Cells.Item(9, 6).Value = 20
$ws1.Cells.Item(11, 6).Value = 510
$ws1.Cells.Item(12, 6).Value = 466
$ws1.Cells.Item(15, 6).Value = 1328
$ws1.Cells.Item(16, 6).Value = 2681
$ws1.Cells.Item(17, 6).Value = 383
$ws1.Cells.Item(18, 6).Value = 83
$ws1.Cells.Item(19, 6).Value = 61
$ws1.Cells.Item(20, 6).Value = 64
$ws1.Cells.Item(21, 6).Value = 2319
$ws1.Cells.Item(22, 6).Value = 94
$ws1.Cells.Item(23, 6).Value = 74
$ws1.Cells.Item(26, 6).Value = 110
$ws1.Cells.Item(28, 6).Value = 228

# Sheet: 演出 (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 40

# Sheet: 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 877
$ws4.Cells.Item(3, 6).Value = 259
$ws4.Cells.Item(4, 6).Value = 20
$ws4.Cells.Item(5, 6).Value = 164
$ws4.Cells.Item(6, 6).Value = 151
$ws4.Cells.Item(7, 6).Value = 157
$ws4.Cells.Item(8, 6).Value = 40
$ws4.Cells.Item(9, 6).Value = 4640
$ws4.Cells.Item(10, 6).Value = 20
$ws4.Cells.Item(12, 6).Value = 510
$ws4.Cells.Item(13, 6).Value = 466
$ws4.Cells.Item(16, 6).Value = 1328
$ws4.Cells.Item(17, 6).Value = 2681
$ws4.Cells.Item(18, 6).Value = 383
$ws4.Cells.Item(19, 6).Value = 83
$ws4.Cells.Item(20, 6).Value = 61
$ws4.Cells.Item(21, 6).Value = 64
$ws4.Cells.Item(22, 6).Value = 2319
$ws4.Cells.Item(23, 6).Value = 94
$ws4.Cells.Item(24, 6).Value = 74
$ws4.Cells.Item(27, 6).Value = 110
$ws4.Cells.Item(29, 6).Value = 228
